$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns B:G), reflecting the rolling window update:
# oldest quarter (former row 2) dropped, all rows shift up, new quarter appended at row 11.
$data = @(
    @(-0.05889625796533703, 0.7809309034430243, 1.486890174830229, 1.219381062191073, 1.25326832255184, 18),
    @(-0.03036361610786248, 0.7674964733547235, 0.9590082094163144, 0.9792896453125165, 1.008943334136136, 17),
    @(0.09993964655644642, 0.6259185880536732, 0.8952003656645542, 0.9461502870393024, 0.9717132266730971, 16),
    @(0.03831129326750842, 0.7119710808408777, 1.080030924484957, 1.039245362984583, 1.074989951428016, 15),
    @(0.09938456961952043, 0.7378794059141224, 1.55852155670565, 1.248407608397854, 1.291421985981127, 14),
    @(0.06981294092732569, 0.5769263726184313, 0.6329722743072022, 0.7955955469377655, 0.8248878488730877, 13),
    @(0.2198501169052137, 0.8740014322694108, 1.535812380474684, 1.239278976048042, 1.273853840583311, 12),
    @(0.1257302919151607, 0.7877651198835861, 0.9791626505880603, 0.9895264779620908, 1.029412453155336, 11),
    @(0.2215020228704644, 0.8099270631829588, 1.114429208763017, 1.055665292014006, 1.08799819374202, 10),
    @(0.2277176327227284, 0.7851486816968395, 0.9325366040153429, 0.9656793484461303, 0.995372515309124, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
